$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")

# Step 1: Insert blank row at 79 (shifts everything below down)
$ws.Rows(79).Insert()

# Step 2: Build correct content+format in a neutral scratch row far below (row 200)
$ws.Range("A200:G200").Value = $ws.Range("A78:G78").Value
$srcFmt = $ws.Range("A78:G78")
$srcFmt.Copy()
$dstFmt = $ws.Range("A200:G200")
$dstFmt.PasteSpecial(-4122)
$dstFmt.Interior.Color = 65535

# Step 3: copy that scratch row into row 79 (values+format)
$ws.Range("A200:G200").Copy()
$ws.Range("A79:G79").PasteSpecial(-4104)

# Step 4: clear scratch row completely
$ws.Range("A200:G200").Clear()
